# "filter headers added, step 2 flu working"
#
# Adds the filter header row + a blank/default filter-value row to the
# top of the first worksheet (H7N3Seq1), then removes the two rows that
# used to sit between that block and the trailing spacer row so the
# spacer (originally row 14) collapses up to row 12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: filter headers
$ws.Range("A1").Value = "H7N3-All"
$ws.Range("B1").Value = "H7N3-01"

# Row 2: default filter values
$ws.Range("A2").Value = "0"
$ws.Range("B2").Value = "0"

# Remove the two now-obsolete rows, shifting everything below (the
# formatted spacer row that used to be row 14) up by two rows
$ws.Range("A3:A4").EntireRow.Delete()

# Touch the spacer row so its column-span metadata is regenerated to
# cover the new A:B data range, keeping it visually/structurally
# unstyled (matches the sheet's default column style)
$ws.Range("A12:B12").Style = "Normal"

# Leave the selection where the user ended up
[void]$ws.Range("E11").Select()
